$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 20, shifting existing rows 20-43 down to 22-45
$ws.Range("A20:A21").EntireRow.Insert()

# Populate new row 20 (Castle Brite / Primera, week of 2022-12-07)
$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "Macroferia Regional de Talca"
$ws.Range("C20").Value = "Maule"
$ws.Range("D20").Value = 44902
$ws.Range("E20").Value = 7
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100103
$ws.Range("H20").Value = "Frutos de hueso (carozo)"
$ws.Range("I20").Value = 100103003
$ws.Range("J20").Value = "Damasco"
$ws.Range("K20").Value = "Castle Brite"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 560
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 15000
$ws.Range("P20").Value = 14643
$ws.Range("Q20").Value = "`$/caja 10 kilos"
$ws.Range("R20").Value = "Región de O'Higgins"
$ws.Range("S20").Value = 1464
$ws.Range("T20").Value = 10

# Populate new row 21 (Castle Brite / Segunda, week of 2022-12-07)
$ws.Range("A21").Value = 5
$ws.Range("B21").Value = "Macroferia Regional de Talca"
$ws.Range("C21").Value = "Maule"
$ws.Range("D21").Value = 44902
$ws.Range("E21").Value = 7
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100103
$ws.Range("H21").Value = "Frutos de hueso (carozo)"
$ws.Range("I21").Value = 100103003
$ws.Range("J21").Value = "Damasco"
$ws.Range("K21").Value = "Castle Brite"
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 190
$ws.Range("N21").Value = 13000
$ws.Range("O21").Value = 13000
$ws.Range("P21").Value = 13000
$ws.Range("Q21").Value = "`$/caja 10 kilos"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 1300
$ws.Range("T21").Value = 10
